# Adds the "CREACION DE AMBIENTES" section (virtual environment setup
# instructions) at the end of the document, right after the last existing
# paragraph ("...libffi-dev python3-dev") and before the final paragraph
# mark / sectPr.

$d = $word.ActiveDocument

# Locate the very end of the document body (just before the last, implicit
# paragraph mark) so the new content is appended as new paragraphs rather
# than replacing anything.
$lastPara = $d.Paragraphs.Last
$insertPoint = $lastPara.Range.End - 1
$target = $d.Range($insertPoint, $insertPoint)

# Raw WordprocessingML for the new paragraphs, inserted verbatim so that
# formatting (bold headings, xml:space="preserve" runs, w:proofErr spell
# check markers, etc.) matches exactly.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>CREACION DE AMBIENTES:</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Verificar donde esta </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>python</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> y </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>pip</w:t>
  </w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>which</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> python3</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>which</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pip3</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Si estas en </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>linus</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> o </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>wsl</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> debes instalar</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">sudo </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>apt</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -y python3-venv</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Poner cada proyecto en su propio ambiente, entrar en cada carpeta</w:t>
  </w:r><w:r><w:t>.</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">python3 -m </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>venv</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>env</w:t>
  </w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Activar el ambiente</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>source</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>env</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bin</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t>
  </w:r><w:r><w:t>actívate</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Salir del ambiente virtual</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>D</w:t>
  </w:r><w:r><w:t>eactivate</w:t>
  </w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Podemos instalar las </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>librerias</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> necesarias en el ambiente virtual </w:t>
  </w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>como</w:t>
  </w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> por ejemplo</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">pip3 </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t>
  </w:r><w:proofErr w:type="spellStart"/><w:r><w:t>matplotlib</w:t>
  </w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>==3.5.0</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Verificar las instalaciones</w:t>
  </w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>pip3 freeze</w:t></w:r></w:p>
'@

$target.InsertXML($xml) | Out-Null
